# Actualizacion tipos de datos BSCS
$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wsMetodos = $wb.Worksheets.Item("Hoja1")
$wsMetodos.Name = "Metodos"

$wsTipos = $wb.Worksheets.Item("Hoja2")
$wsTipos.Name = "Tipos de Datos"

# --- Fix a value on "Metodos" (C3: OM -> OE) ---
$wsMetodos.Range("C3").Value = "OE"

# --- Move selection on "Metodos" sheet to C4 ---
[void]$wsMetodos.Range("C4").Select()

# --- Populate "Tipos de Datos" sheet with new BSCS data-type table ---
$wsTipos.Range("A1").Value = "Nombre"
$wsTipos.Range("B1").Value = "Tipo"
$wsTipos.Range("C1").Value = "Largo"
$wsTipos.Range("D1").Value = "Descripcion"

$wsTipos.Range("A2").Value = "ESTADO"
$wsTipos.Range("B2").Value = "STRING"
$wsTipos.Range("C2").Value = 32
$wsTipos.Range("D2").Value = "Estado del Sistema"

$wsTipos.Range("A3").Value = "ID_CLIENTE"
$wsTipos.Range("B3").Value = "STRING"
$wsTipos.Range("C3").Value = 16
$wsTipos.Range("D3").Value = "Identificador unico del Cliente. Llave primaria de la tabla Clientes"

$wsTipos.Range("A4").Value = "IMSI"
$wsTipos.Range("B4").Value = "INTEGER"
$wsTipos.Range("C4").Value = 15

$wsTipos.Range("A5").Value = "MSISDN"
$wsTipos.Range("B5").Value = "INTEGER"
$wsTipos.Range("C5").Value = 15

$wsTipos.Range("D4").Value = "Imsi del abonado"
$wsTipos.Range("D5").Value = "IMSISDN del Abonado"

$wsTipos.Range("A6").Value = "PLAN_BSCS"
$wsTipos.Range("B6").Value = "STRING"
$wsTipos.Range("C6").Value = 32
$wsTipos.Range("D6").Value = "Plan En BSCS"

# --- Column D width (best-fit for the long description text) ---
$wsTipos.Columns.Item(4).ColumnWidth = 59

# --- Make "Tipos de Datos" the active/visible tab with its own selection ---
[void]$wsTipos.Range("A7").Select()
$wsTipos.Activate()

Write-Host "done"
